# "Fruta / hortaliza, semanal" — weekly refresh of the Zanahoria
# (Terminal Hortofrutícola Agro Chillán) subset: a new week's record is
# inserted above the existing data at row 468, shifting every subsequent
# row down by one (dimension grows from A1:R511 to A1:R512).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 468 — pushes old rows 468..511 down
# to 469..512, carrying their values/styles with them.
$ws.Rows.Item(468).Insert()

# Populate the newly inserted row with the latest weekly observation.
$ws.Cells.Item(468, 1).Value2  = 7
$ws.Cells.Item(468, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(468, 3).Value2  = "Ñuble"
$ws.Cells.Item(468, 4).Value2  = 45166
$ws.Cells.Item(468, 5).Value2  = 16
$ws.Cells.Item(468, 6).Value2  = 100114013
$ws.Cells.Item(468, 7).Value2  = "Zanahoria"
$ws.Cells.Item(468, 8).Value2  = "Sin especificar"
$ws.Cells.Item(468, 9).Value2  = "Primera"
$ws.Cells.Item(468, 10).Value2 = 250
$ws.Cells.Item(468, 11).Value2 = 6000
$ws.Cells.Item(468, 12).Value2 = 7000
$ws.Cells.Item(468, 13).Value2 = 6600
$ws.Cells.Item(468, 14).Value2 = "$/saco 20 kilos"
$ws.Cells.Item(468, 15).Value2 = "Provincia de Diguillín"
$ws.Cells.Item(468, 16).Value2 = 330
$ws.Cells.Item(468, 17).Value2 = 20
$ws.Cells.Item(468, 18).Value2 = "Hortaliza"
